$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Fix the address locator text: "US Highway" -> "U.S. Highway" (and
# re-order of the shared-string table that Excel performs as a result).
# The leading apostrophe preserves the original "quote prefix" / forced-text
# formatting these two cells already had, instead of Excel clearing it.
$ws.Range("F2").Value = "'2884 U.S. Highway 67, Farmington, MO 63640"
$ws.Range("G2").Value = "'2884 US-67"

# Scroll the view so column C is left-most, then move the active selection
# to F7, matching the saved view/selection state.
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 1
$ws.Range("F7").Select()
